# Actualización al 27 de octubre 2023
$wb = $excel.ActiveWorkbook

# --- Sheet "Ingreso": update last row (544) and append new aportes rows ---
$ingreso = $wb.Worksheets.Item("Ingreso")

$ingresoRows = @(
    @(45214, "Robert", 400, "Aporte"),
    @(45214, "Kibelo", 1000, "Aporte"),
    @(45214, "Wilkin", 100, "Aporte"),
    @(45214, "Kukito", 50, "Aporte"),
    @(45214, "Jeicol", 100, "Aporte"),
    @(45214, "Melvin", 100, "Aporte"),
    @(45214, "Omaury", 100, "Aporte"),
    @(45214, "Rayder", 100, "Aporte"),
    @(45214, "Yandi", 300, "Aporte"),
    @(45214, "Anuel", 200, "Aporte")
)

$startRow = 544
for ($i = 0; $i -lt $ingresoRows.Count; $i++) {
    $r = $startRow + $i
    $data = $ingresoRows[$i]
    $ingreso.Cells.Item($r, 1).Value2 = $data[0]
    $ingreso.Cells.Item($r, 2).Value2 = $data[1]
    $ingreso.Cells.Item($r, 3).Value2 = $data[2]
    $ingreso.Cells.Item($r, 4).Value2 = $data[3]
}

# --- Sheet "Gastos": update row 64 amount and append new row 66 ---
$gastos = $wb.Worksheets.Item("Gastos")

$gastos.Range("C64").Value2 = 1200

$gastos.Cells.Item(66, 1).Value2 = 45214
$gastos.Cells.Item(66, 2).Value2 = "Agua y arbitro"
$gastos.Cells.Item(66, 3).Formula = "=400+230"

# Leave the view positioned on the cell that was last edited on this sheet.
[void]$gastos.Range("C64").Select()

# Re-activate "Ingreso" last so it remains the selected/visible tab, with the
# cursor parked on the next empty entry row, matching the saved view state.
[void]$ingreso.Activate()
[void]$ingreso.Range("B554").Select()
